$d = $word.ActiveDocument

# "> Create user stories (1 hour)"  ->  "> Create user stories (15 minutes)"
$d.Content.Find.Execute("Create user stories (1 hour)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Create user stories (15 minutes)", 2)

# "> Finish work on deliverable 3 (1 hour)"  ->  "> Finish work on deliverable 3 (10 minutes)"
$d.Content.Find.Execute("Finish work on deliverable 3 (1 hour)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Finish work on deliverable 3 (10 minutes)", 2)
